# Update previous_count (column C) and change (column D) for agencies
# whose counts now match the current_count (column B), reflecting the
# new archived reference file used for the comparison.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Department of Commerce
$ws.Cells.Item(6, 3).Value = 27
$ws.Cells.Item(6, 4).Value = 0

# Row 10: Department of Health and Human Services
$ws.Cells.Item(10, 3).Value = 45
$ws.Cells.Item(10, 4).Value = 0

# Row 11: Department of Homeland Security
$ws.Cells.Item(11, 3).Value = 38
$ws.Cells.Item(11, 4).Value = 0

# Row 13: Department of Justice
$ws.Cells.Item(13, 3).Value = 32
$ws.Cells.Item(13, 4).Value = 0

# Row 21: Executive Office of the President, Management and Administration
$ws.Cells.Item(21, 3).Value = 6
$ws.Cells.Item(21, 4).Value = 0

# Row 29: National Security Council
$ws.Cells.Item(29, 3).Value = 19
$ws.Cells.Item(29, 4).Value = 0

# Row 32: Office of Personnel Management
$ws.Cells.Item(32, 3).Value = 11
$ws.Cells.Item(32, 4).Value = 0

# Row 35: Small Business Administration
$ws.Cells.Item(35, 3).Value = 19
$ws.Cells.Item(35, 4).Value = 0

# Row 40: United States Mission to the United Nations
$ws.Cells.Item(40, 3).Value = 10
$ws.Cells.Item(40, 4).Value = 0

# Row 41: United States Postal Service
$ws.Cells.Item(41, 3).Value = 5
$ws.Cells.Item(41, 4).Value = 0
